# Update handback-status workbook: replace the two file identities with new ones
# and refresh the generated/handback timestamps across the Overview, zh-cn and
# de-de sheets (this mirrors a fresh "Generate Report for Handback" run).

$wb = $excel.ActiveWorkbook

$oldId1 = "5b07f70b-9de1-4c36-9ce1-80d7ea4c5ade"
$newId1 = "627c0b49-119a-44fb-abba-4b27c494d0fa"
$oldId2 = "7a13037e-e344-4ba5-8413-f3aabe8bfefc"
$newId2 = "ffff2a12ba4d-51fd-42fa-a0f5-bf502d17cf49"

$newXlf  = "627c0b49-119a-44fb-abba-4b27c494d0fa.015b6d9c86507b0cb106e95d72d34fb79f019acc"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId1.md"
$wsOverview.Range("B2").Value = "e2e\$newId1.md"
$wsOverview.Range("G2").Value = "2016-08-18 21:03:54"

$wsOverview.Range("A3").Value = "$newId2.md"
$wsOverview.Range("B3").Value = "e2e\$newId2.md"
$wsOverview.Range("G3").Value = "2016-08-18 21:03:54"

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.TextToDisplay -eq "e2e\$oldId1.md") {
        $h.TextToDisplay = "e2e\$newId1.md"
    } elseif ($h.TextToDisplay -eq "e2e\$oldId2.md") {
        $h.TextToDisplay = "e2e\$newId2.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId1.md"
$wsZhCn.Range("G2").Value = "$newXlf.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-18 21:03:49"
$wsZhCn.Range("I2").Value = "$newId1.md"
$wsZhCn.Range("J2").Value = "$newXlf.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-18 21:04:14"

$wsZhCn.Range("A3").Value = "$newId2.md"
$wsZhCn.Range("G3").Value = "$newXlf.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-18 21:03:49"
$wsZhCn.Range("I3").Value = "$newId2.md"
$wsZhCn.Range("J3").Value = "$newXlf.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-18 21:04:14"

foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.TextToDisplay -eq "$oldId1.md") {
        $h.TextToDisplay = "$newId1.md"
    } elseif ($h.TextToDisplay -eq "$oldId2.md") {
        $h.TextToDisplay = "$newId2.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId1.md"
$wsDeDe.Range("G2").Value = "$newXlf.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-18 21:03:54"
$wsDeDe.Range("I2").Value = "$newId1.md"
$wsDeDe.Range("J2").Value = "$newXlf.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-18 21:04:22"

$wsDeDe.Range("A3").Value = "$newId2.md"
$wsDeDe.Range("G3").Value = "$newXlf.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-18 21:03:54"
$wsDeDe.Range("I3").Value = "$newId2.md"
$wsDeDe.Range("J3").Value = "$newXlf.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-18 21:04:22"

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.TextToDisplay -eq "$oldId1.md") {
        $h.TextToDisplay = "$newId1.md"
    } elseif ($h.TextToDisplay -eq "$oldId2.md") {
        $h.TextToDisplay = "$newId2.md"
    }
}
